$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 432-433; everything currently at 432.. shifts down by 2.
$ws.Rows("432:433").Insert()

# New row 432: Hass / Primera, Cabildo, caja de 17 kilos
$ws.Cells.Item(432, 1).Value = 11
$ws.Cells.Item(432, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(432, 3).Value = "Bíobío"
$ws.Cells.Item(432, 4).Value = 44595
$ws.Cells.Item(432, 5).Value = 8
$ws.Cells.Item(432, 6).Value = "Fruta"
$ws.Cells.Item(432, 7).Value = 100106
$ws.Cells.Item(432, 8).Value = "Oleaginosos"
$ws.Cells.Item(432, 9).Value = 100106002
$ws.Cells.Item(432, 10).Value = "Palta"
$ws.Cells.Item(432, 11).Value = "Hass"
$ws.Cells.Item(432, 12).Value = "Primera"
$ws.Cells.Item(432, 13).Value = 270
$ws.Cells.Item(432, 14).Value = 2500
$ws.Cells.Item(432, 15).Value = 2800
$ws.Cells.Item(432, 16).Value = 2667
$ws.Cells.Item(432, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(432, 18).Value = "Cabildo"
$ws.Cells.Item(432, 19).Value = 2667
$ws.Cells.Item(432, 20).Value = 1

# New row 433: Hass / Segunda, Cabildo, caja de 17 kilos
$ws.Cells.Item(433, 1).Value = 11
$ws.Cells.Item(433, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(433, 3).Value = "Bíobío"
$ws.Cells.Item(433, 4).Value = 44595
$ws.Cells.Item(433, 5).Value = 8
$ws.Cells.Item(433, 6).Value = "Fruta"
$ws.Cells.Item(433, 7).Value = 100106
$ws.Cells.Item(433, 8).Value = "Oleaginosos"
$ws.Cells.Item(433, 9).Value = 100106002
$ws.Cells.Item(433, 10).Value = "Palta"
$ws.Cells.Item(433, 11).Value = "Hass"
$ws.Cells.Item(433, 12).Value = "Segunda"
$ws.Cells.Item(433, 13).Value = 200
$ws.Cells.Item(433, 14).Value = 2000
$ws.Cells.Item(433, 15).Value = 2000
$ws.Cells.Item(433, 16).Value = 2000
$ws.Cells.Item(433, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(433, 18).Value = "Cabildo"
$ws.Cells.Item(433, 19).Value = 2000
$ws.Cells.Item(433, 20).Value = 1
